# SSDM-12286 Fixed letter case inconsistencies.
# "Vocabulary Code" -> "Vocabulary code" and "Generated Code Prefix" -> "Generated code prefix"
# in the header rows of both tables on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 4 / row 11: "Vocabulary Code" -> "Vocabulary code"
$ws.Range("H4").Value = "Vocabulary code"
$ws.Range("H11").Value = "Vocabulary code"

# Header row 2 / row 9: "Generated Code Prefix" -> "Generated code prefix"
$ws.Range("E2").Value = "Generated code prefix"
$ws.Range("E9").Value = "Generated code prefix"

# Matches the editor's final cursor position recorded in the saved file.
$ws.Range("E9").Select()
